# Commit: "Separated DateTime so that Date and Time are separate strings
# allowing for more flexibility."
#
# The underlying change recorded in the diff is the addition of a new batch
# of timestamp log entries to the shared-string table, with the sheet's
# tracked/most-recent timestamp cell (B2) moved forward to the newest one
# ("03/30/2020 23:57:28" - the last of the newly appended timestamps).
#
# Update the tracked timestamp cell to the newest logged value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "03/30/2020 23:57:28"
